$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-19 from
# 2023-09-14 (serial 45183) to 2023-09-15 (serial 45184).
$ws.Range("C2:C19").Value = 45184
